$wb = $excel.ActiveWorkbook

$ws3 = $wb.Worksheets.Item("Uniform - Random")

# --- 1. Duplicate "Uniform - Random" to the end of the workbook, then rename it.
#        Doing this before touching $ws3's own values means the new sheet starts
#        out as a faithful clone (same shared strings, column widths, formulas).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3.Copy([System.Type]::Missing, $lastSheet)
$ws4 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4.Name = "Proportional - Constant"

# --- 2. Rewrite the new "Proportional - Constant" sheet's data ---
$ws4.Range("C2").Value = 59893
$ws4.Range("D2").Value = 0
$ws4.Range("E2").Value = 0

$ws4.Range("C3").Value = 0
$ws4.Range("D3").Value = 0
$ws4.Range("E3").Value = 0

$ws4.Range("C4").Value = 0
$ws4.Range("D4").Value = 0
$ws4.Range("E4").Value = 0

$ws4.Range("C5").Value = 0
$ws4.Range("D5").Value = 0
$ws4.Range("E5").Value = 0

$ws4.Range("C6").Value = 0
$ws4.Range("D6").Value = 0
$ws4.Range("E6").Value = 0

$ws4.Range("C7").Value = 0
$ws4.Range("D7").Value = 0
$ws4.Range("E7").Value = 0

$ws4.Range("C8").Value = 0
$ws4.Range("D8").Value = 0
$ws4.Range("E8").Value = 0

$ws4.Range("C9").Value = 0
$ws4.Range("D9").Value = 0
$ws4.Range("E9").Value = 0

$ws4.Range("C10").Value = 0
$ws4.Range("D10").Value = 0
$ws4.Range("E10").Value = 0

$ws4.Range("C11").Value = 0
$ws4.Range("D11").Value = 0
$ws4.Range("E11").Value = 0

# D12/E12 keep their AVERAGE formulas (copied from the source sheet); they will
# recalc automatically to 0 once the rows above all read 0.

# New sheet's selection: single cell C2 is selected/active.
$ws4.Range("C2").Select() | Out-Null

# --- 3. Update the existing "Uniform - Random" sheet's own numbers ---
$ws3.Range("D8").Value = 29516
$ws3.Range("E8").Value = 75761

$ws3.Range("C9").Value = 44656
$ws3.Range("D9").Value = 29556
$ws3.Range("E9").Value = 76138

$ws3.Range("C10").Value = 44695
$ws3.Range("D10").Value = 29402
$ws3.Range("E10").Value = 75138

$ws3.Range("C11").Value = 44689
$ws3.Range("D11").Value = 29331
$ws3.Range("E11").Value = 75286

# "Uniform - Random" is no longer the active tab; its selection becomes the
# whole used range instead of the old single-cell D8 selection. (Selecting on
# it also makes it the active sheet for a moment, hence the re-activation of
# "Proportional - Constant" below.)
$ws3.Range("A1:E12").Select() | Out-Null

# Leave "Proportional - Constant" as the active sheet/tab (matches activeTab=3)
# with its C2 selection intact.
$ws4.Activate() | Out-Null
